$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new ALU op "A++" ---
$ws.Range("A10").Value = "x"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = "A++"

# --- Row 11: new ALU op "A--" ---
$ws.Range("A11").Value = "x"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = "A--"

# Replicate the formatting of the previous data row (row 9) onto the two
# new rows so they keep the same table styling (borders/fills) used by
# the rest of the data rows.
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H9").Copy()
$ws.Range("H10:H11").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Match the workbook's recorded selection after the edit.
$ws.Range("D11").Select()
